# Regenerate Tests Results: replace search-time/iteration data in columns A and B
# and adjust the "average/100" formulas in E4/F4. Also update the active
# selection to I6 (matches the author's post-edit cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New sample data (rows 2..51) for column A ("Время работы алгоритмы")
$colA = @(
    134100, 284700, 43700, 16500, 14200,
    8300, 8800, 7700, 9300, 9700,
    9500, 9500, 9800, 6400, 15300,
    7900, 7400, 6700, 10100, 7000,
    5600, 6500, 10300, 9400, 10900,
    8100, 8900, 9700, 7500, 6000,
    10400, 7500, 11500, 6700, 6900,
    11100, 10200, 9700, 12000, 6400,
    6100, 7100, 8500, 8800, 8400,
    8200, 9700, 8900, 8400, 6900
)

# New sample data (rows 2..51) for column B ("Количество итераций")
$colB = @(
    1900, 1900, 1900, 1900, 1900,
    1900, 1900, 1900, 1900, 1900,
    1900, 1900, 1900, 1900, 1900,
    1900, 1900, 1900, 1900, 1900,
    1900, 1900, 1900, 1900, 1900,
    1900, 1900, 1900, 1900, 1900,
    1900, 1900, 1900, 1900, 1900,
    1900, 1900, 1900, 1900, 1900,
    1900, 1900, 1900, 1900, 1900,
    1900, 1900, 1900, 1900, 1900
)

for ($i = 0; $i -lt $colA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# Update summary formulas: was /10000, now /100
$ws.Range("E4").Formula = "=E2/100"
$ws.Range("F4").Formula = "=F2/100"

# Update the selected cell shown when the sheet was last saved
$ws.Range("I6").Select()

$wb.Save()
